$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the time-types cell for the average-blood-pressure component row: "Periodĵ" -> "dateTimeĵ, Periodĵ"
$ws.Range("G2").Value = "dateTimeĵ, Periodĵ"

# Row 6 previously held the "us-core-treatment-intervention-preference" profile row (with a
# data bug: the "Code" column duplicated the Category Code text). It should instead hold the
# "us-core-observation-lab" profile row.
$ws.Range("A6").Value = "us-core-observation-lab"
$ws.Range("B6").Value = "US Core Laboratory Result Observation Profile"
$ws.Range("C6").Value = "Observation Category Codes#laboratory"
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "http://hl7.org/fhir/us/core/ValueSet/us-core-laboratory-test-codes (extensible)"

# Row 7 previously held the "us-core-observation-lab" profile row. It should instead hold the
# "us-core-treatment-intervention-preference" profile row, with the Category Code text corrected
# to "US Core Category#treatment-intervention-preference" (also duplicated into "Code", matching
# the original data bug pattern).
$ws.Range("A7").Value = "us-core-treatment-intervention-preference"
$ws.Range("B7").Value = "US Core Treatment Intervention Preference Profile"
$ws.Range("C7").Value = "US Core Category#treatment-intervention-preference"
$ws.Range("E7").Value = "US Core Category#treatment-intervention-preference"
$ws.Range("F7").Value = ""
